$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Periodo Mora" column (E) to ascending order 1907 -> 2001
$ws.Range("E16").Value = "1907"
$ws.Range("E17").Value = "1908"
$ws.Range("E18").Value = "1909"
$ws.Range("E19").Value = "1910"
$ws.Range("E20").Value = "1911"
$ws.Range("E21").Value = "1912"
$ws.Range("E22").Value = "2001"

# Update "Valor Mora" column (F)
$ws.Range("F16").Value = 51600
$ws.Range("F17").Value = 51600
$ws.Range("F18").Value = 51600
$ws.Range("F19").Value = 51600
$ws.Range("F20").Value = 51600
$ws.Range("F21").Value = 51600
$ws.Range("F22").Value = 36120

# Update "Salario Basico" column (G)
$ws.Range("G16").Value = 1290000
$ws.Range("G17").Value = 1290000
$ws.Range("G18").Value = 1290000
$ws.Range("G19").Value = 1290000
$ws.Range("G20").Value = 1290000
$ws.Range("G21").Value = 1290000
$ws.Range("G22").Value = 1290000
